$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the row40/row41 coin swap)
# Every assignment is apostrophe-prefixed to force text interpretation (these
# columns store numbers/percentages as plain text, e.g. "58.040.86", "1.00"),
# then the cell style is reset to "Normal" so no stray quote-prefix style sticks.
$ws.Range("D2").Value = "'58.040.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.451.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.47%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'524.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.02%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'130.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.22%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.456.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.11%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.324"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.888.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'58.019.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.456.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.10%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.52%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'315.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.38%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.14%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.67%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.567.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.56%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.157"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.68%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'173.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.29%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0737"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.11%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.80%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.21%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'17.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.24%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -5.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.65%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'OKB"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'36.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.49%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'SuiNetwork"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.53%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.69%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.80%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'261.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.20%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.585"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.21%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.59%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0921"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'122.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.02%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0493"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'17.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.00%  "
$ws.Range("E51").Style = "Normal"
